$d = $word.ActiveDocument

$d.Content.Find.Execute("76÷7=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "47÷4=11, 3", 2) | Out-Null
$d.Content.Find.Execute("29÷9=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷5=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "57÷5=11, 2", 2) | Out-Null
$d.Content.Find.Execute("31÷8=3, 7", $true, $false, $false, $false, $false, $true, 1, $false, "23÷5=4, 3", 2) | Out-Null
$d.Content.Find.Execute("31÷9=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=47, 0", 2) | Out-Null
$d.Content.Find.Execute("49÷4=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=5, 1", 2) | Out-Null
$d.Content.Find.Execute("34÷7=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "82÷7=11, 5", 2) | Out-Null
$d.Content.Find.Execute("22÷3=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "95÷6=15, 5", 2) | Out-Null
$d.Content.Find.Execute("84÷7=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷9=9, 1", 2) | Out-Null
$d.Content.Find.Execute("54÷4=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "20÷9=2, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷2=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=25, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷4=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "48÷9=5, 3", 2) | Out-Null
$d.Content.Find.Execute("86÷7=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "27÷3=9, 0", 2) | Out-Null
$d.Content.Find.Execute("30÷2=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷2=48, 0", 2) | Out-Null
$d.Content.Find.Execute("69÷9=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=4, 0", 2) | Out-Null
$d.Content.Find.Execute("86÷4=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 2) | Out-Null
$d.Content.Find.Execute("75÷9=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "76÷7=10, 6", 2) | Out-Null
$d.Content.Find.Execute("37÷6=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "28÷3=9, 1", 2) | Out-Null
$d.Content.Find.Execute("96÷7=13, 5", $true, $false, $false, $false, $false, $true, 1, $false, "79÷6=13, 1", 2) | Out-Null
$d.Content.Find.Execute("68÷4=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2) | Out-Null
$d.Content.Find.Execute("92÷3=30, 2", $true, $false, $false, $false, $false, $true, 1, $false, "96÷5=19, 1", 2) | Out-Null
$d.Content.Find.Execute("99÷3=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "83÷6=13, 5", 2) | Out-Null
$d.Content.Find.Execute("78÷4=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "60÷5=12, 0", 2) | Out-Null
$d.Content.Find.Execute("19÷6=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=18, 3", 2) | Out-Null
$d.Content.Find.Execute("56÷3=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "39÷2=19, 1", 2) | Out-Null
